$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update state names in column A (rows 2-26) to properly capitalized / spaced values.
$names = @(
    "Andhra Pradesh",
    "Arunachal Pradesh",
    "Assam",
    "Bihar",
    "Delhi",
    "Goa",
    "Gujarat",
    "Haryana",
    "Himachal Pradesh",
    "Jammu & Kashmir",
    "Karnatka",
    "Kerala",
    "Madhya Pradesh",
    "Maharashtra",
    "Manipur",
    "Meghalya",
    "Mizoram",
    "Nagaland",
    "Odisha",
    "Punjab",
    "Rajasthan",
    "Tamil Nadu",
    "Tripura",
    "Uttar Pradesh",
    "West Bengal"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
}

# Widen column A slightly (matches the author's manual resize).
$ws.Columns.Item(1).ColumnWidth = 13.92

# Update the active view: zoom to 135% and move the selection to B20.
$excel.ActiveWindow.Zoom = 135
$ws.Range("B20").Select()
